$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Append four new feedback rows (57-60) to the bottom of the log, mirroring
# the formatting of the existing rows (date columns A/F use the m/d/yyyy
# style from column A56, description/status/action columns D/E/H use the
# wrap-text style from column D56/E56).
# ---------------------------------------------------------------------------

function Set-LikeA($srcAddr, $dstAddr, $val) {
    $ws.Range($srcAddr).Copy($ws.Range($dstAddr)) | Out-Null
    $ws.Range($dstAddr).Value = $val
}

# Pre-stage the date/text formatting (style only, no new shared strings) for
# every new cell first ...
Set-LikeA "A56" "A57" 41935
Set-LikeA "F56" "F57" 41935
Set-LikeA "A56" "A58" 41933
Set-LikeA "F56" "F58" 41933
Set-LikeA "A56" "A59" 41933
Set-LikeA "F56" "F59" 41933
Set-LikeA "A56" "A60" 41933
Set-LikeA "F56" "F60" 41933

$ws.Range("B57").Value = "Email"
$ws.Range("C57").Value = "John Yang <john.y@solaronesolution.com>"
$ws.Range("B58").Value = "Email"
$ws.Range("C58").Value = "John Yang <john.y@solaronesolution.com>"
$ws.Range("B59").Value = "Email"
$ws.Range("B60").Value = "Email"

$ws.Range("D56").Copy($ws.Range("D57")) | Out-Null
$ws.Range("D56").Copy($ws.Range("D58")) | Out-Null
$ws.Range("D56").Copy($ws.Range("H58")) | Out-Null
$ws.Range("D56").Copy($ws.Range("D59")) | Out-Null
$ws.Range("D56").Copy($ws.Range("D60")) | Out-Null
$ws.Range("E56").Copy($ws.Range("E57")) | Out-Null
$ws.Range("E56").Copy($ws.Range("E58")) | Out-Null
$ws.Range("E56").Copy($ws.Range("E59")) | Out-Null
$ws.Range("E56").Copy($ws.Range("E60")) | Out-Null

# ... then assign the actual text values in the exact order the strings were
# first introduced, so the shared-string table comes out in the same order
# as the author's own edit (new text landed in the pool as it was typed:
# row 57 description, row 58 description, row 58 action, row 59 description,
# row 59 contact, row 57 status, row 58 status, row 60 description, row 60
# contact).
$ws.Range("D57").Value = "Here is the model I would like to simulate.  The model has 208 modules on the flat roof and  70 on the pitched roof.  I modeled this building looks like two but actually it is just one building; so, just one electricity meter.  At least two different string inverters are needed.`nI believe that, for this situation, multiple subsystem modeling feature is much useful.`n"
$ws.Range("D58").Value = ".  I am attaching two screen shots; one is for SAM 2014.1.14  and another is for SAM 2014.9.30.  For PG&E residential, usually it is monthly tier but on new SAM 2014.9.30, there is no monthly tiers."
$ws.Range("H58").Value = "Usability issue with URDB window"
$ws.Range("D59").Value = "The ‘register’ button didn’t look like a button. "
$ws.Range("C59").Value = "Michael F. Troge <mtroge@oneidanation.org>"
$ws.Range("E57").Value = "Followed up. Forwarded to team."
$ws.Range("E58").Value = "Followed up. On meeting agenda."
$ws.Range("D60").Value = "Will  the new version be able to load saved zsam files from the previous version, or TMY.tm2 weather files downloaded in the previous version?"
$ws.Range("C60").Value = "Gomez, Tommaso <tommaso.gomez@intel.com>"

# E59/E60 reuse the pre-existing "Followed up" shared string (index 66).
$ws.Range("E59").Value = "Followed up"
$ws.Range("E60").Value = "Followed up"

$ws.Rows.Item(57).RowHeight = 120
$ws.Rows.Item(58).RowHeight = 45
# Row 59 keeps the sheet's default row height (15) - no explicit ht= in the
# target XML, so leave it untouched.
$ws.Rows.Item(60).RowHeight = 30

# Move the view/selection down to where the new last row is, matching the
# author's saved cursor position when they finished editing.
$ws.Range("A61").Select()
